$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update 想去人数 (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1008
$ws1.Range("F3").Value = 2072
$ws1.Range("F4").Value = 456

# Sheet "全部类型" (All Types) - same three exhibitions appear again, further down
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1008
$ws4.Range("F5").Value = 2072
$ws4.Range("F6").Value = 456
